$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (Agus Budi / Lamongan) ---
$ws.Cells.Item(2,1).Value = 14241010994
$ws.Cells.Item(2,2).Value = "Agus Budi"
$ws.Cells.Item(2,3).Value = "Lamongan"
$ws.Cells.Item(2,4).Value = 81777888999

# --- Update row 3 (Aurel Herman / Jember) ---
$ws.Cells.Item(3,1).Value = 14241010123
$ws.Cells.Item(3,2).Value = "Aurel Herman"
$ws.Cells.Item(3,3).Value = "Jember"
$ws.Cells.Item(3,4).Value = 81625827019

# --- Remove old row 4 (Anriko Chiesa) entirely ---
$ws.Rows.Item(4).Delete()

# --- Column C: narrower width, drop the bestFit auto-sizing flag ---
$ws.Columns.Item(3).ColumnWidth = 19.833333333333336

# --- Column E formatting leftover is no longer needed: shift it out and
#     drop it cleanly (Insert pushes the old col E width descriptor to F,
#     then deleting E:F removes both the blank gap and the old descriptor
#     without leaving a stray/empty <col> entry behind). ---
$ws.Columns.Item(5).Insert()
$ws.Range("E1:F1").EntireColumn.Delete()

# --- Update the active selection shown when the sheet is opened ---
$ws.Range("D10").Select()

Write-Host ("Final UsedRange: " + $ws.UsedRange.Address())
